# Task1: add the "Fourth stage" (gradient boosting) paragraph at the very
# end of the document, matching the formatting of the existing body
# paragraphs (FreeSans, 14pt / sz 28, tab-indented first line).

$d = $word.ActiveDocument

# --- 1. Move to the end of the document and start a brand new paragraph ---
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

$newPara = $d.Paragraphs($d.Paragraphs.Count)
$bodyRange = $newPara.Range
$insertStart = $bodyRange.Start
$paraEnd = $bodyRange.End

# --- 2. Type the paragraph text (it inherits the same run formatting used
#        throughout the rest of the report: FreeSans, 14pt, not bold/italic) ---
$textRange = $d.Range($insertStart, $paraEnd - 1)
$textRange.Text = "Четвертый этап — прогнозирование с помощью алгоритма градиентного бустинга. Будем считать, что текущее значение данного мне ряда зависит от 15 предыдущих. Основываясь на этом, составим матрицу X и столбец y, где строка матрицы — 15 предыдущих значений, а соотвествующий строке элемент столбца — зависимое значение (концепция матрицы признаков и столбца ответов)."

# --- 3. Prefix the paragraph with a tab character (leading indent), as in
#        every other paragraph of this report ---
$tabPos = $d.Range($insertStart, $insertStart)
$tabPos.InsertXML('<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="FreeSans" w:hAnsi="FreeSans"/><w:b w:val="0"/><w:bCs w:val="0"/><w:i w:val="0"/><w:iCs w:val="0"/><w:position w:val="0"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="none"/><w:vertAlign w:val="baseline"/></w:rPr><w:tab/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
